# Glassdoor job data refresh: a prior scrape run lost its "random sleep
# 10-15" delay and got throttled, so several job rows came back as
# "No Title Available" / "No Company Name Available" placeholders. This
# re-run fixes up the titles/companies for a handful of rows that now
# resolved, and appends the newly scraped rows (33-62) at the bottom.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix up rows whose title/company resolved on this pass -------------
$ws.Cells.Item(2, 1).Value = "IT Helpdesk Analyst"
$ws.Cells.Item(2, 2).Value = "Glitz Info Solutions"

$ws.Cells.Item(3, 1).Value = "Sr. IT Recruiter"
$ws.Cells.Item(3, 2).Value = "Mesolith Software"

$ws.Cells.Item(4, 1).Value = "IT Support Specialist"
$ws.Cells.Item(4, 2).Value = "CyntraLabs TechLabs"

$ws.Cells.Item(6, 1).Value = "IT Intern"
$ws.Cells.Item(6, 2).Value = "Pixel Vision Technologies"

$ws.Cells.Item(7, 1).Value = "Solutions Architect"
$ws.Cells.Item(7, 2).Value = "AbroadWorks Inc."

# --- Append newly scraped rows 33-62 ------------------------------------
# Most came back throttled (placeholder text); a few in the 49-53 block
# resolved with real title/company values.
$titles = @{
    49 = "IT Recruiter"
    50 = "IT Executive"
    51 = "IT Systems Administrator"
    52 = "Executive IT"
    53 = "IT - Information Technology"
}
$companies = @{
    49 = "WebPariwar"
    50 = "Pinnacle Infotech"
    51 = "Aastitva being foundation"
    52 = "Lotus Petal Charitable Foundation"
    53 = "Larsen & Toubro"
}

for ($r = 33; $r -le 62; $r++) {
    if ($titles.ContainsKey($r)) {
        $title = $titles[$r]
        $company = $companies[$r]
    } else {
        $title = "No Title Available"
        $company = "No Company Name Available"
    }

    $ws.Cells.Item($r, 1).Value = $title
    $ws.Cells.Item($r, 2).Value = $company
    $ws.Cells.Item($r, 3).Value = "No Location Available"
    $ws.Cells.Item($r, 4).Value = "Not Disclosed"
    $ws.Cells.Item($r, 5).Value = "No Description Available"
    $ws.Cells.Item($r, 6).Value = "No Description Available"
    $ws.Cells.Item($r, 7).Value = "No Rating"
}
